$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue "D2" '64.799.27'
Set-TextValue "E2" '  +2.17%  '
Set-TextValue "D3" '3.462.32'
Set-TextValue "E3" '  +2.05%  '
Set-TextValue "E4" '  -0.06%  '
Set-TextValue "D5" '576.63'
Set-TextValue "E5" '  -0.09%  '
Set-TextValue "D6" '161.91'
Set-TextValue "E6" '  +4.09%  '
Set-TextValue "D7" '0.999'
Set-TextValue "E7" '  +0.02%  '
Set-TextValue "D8" '3.462.71'
Set-TextValue "E8" '  +1.70%  '
Set-TextValue "E9" '  +8.16%  '
Set-TextValue "E10" '  -2.44%  '
Set-TextValue "E11" '  +4.17%  '
Set-TextValue "D12" '0.439'
Set-TextValue "E12" '  +0.66%  '
Set-TextValue "D13" '4.056.62'
Set-TextValue "E13" '  +1.85%  '
Set-TextValue "E14" '  -2.92%  '
Set-TextValue "E15" '  +5.57%  '
Set-TextValue "D16" '28.80'
Set-TextValue "E16" '  +6.42%  '
Set-TextValue "D17" '64.764.22'
Set-TextValue "E17" '  +1.93%  '
Set-TextValue "D18" '3.498.98'
Set-TextValue "E18" '  +2.96%  '
Set-TextValue "D19" '6.37'
Set-TextValue "E19" '  +0.07%  '
Set-TextValue "E20" '  +2.12%  '
Set-TextValue "D21" '389.96'
Set-TextValue "E21" '  +0.90%  '
Set-TextValue "D22" '8.20'
Set-TextValue "E22" '  -3.33%  '
Set-TextValue "D23" '0.548'
Set-TextValue "E23" '  +2.60%  '
Set-TextValue "D24" '73.04'
Set-TextValue "E24" '  +2.99%  '
Set-TextValue "D25" '1.00'
Set-TextValue "E25" '  +0.11%  '
Set-TextValue "D26" '0.0000124'
Set-TextValue "E26" '  +20.12%  '
Set-TextValue "E27" '  -0.91%  '
Set-TextValue "E28" '  +0.61%  '
Set-TextValue "D29" '1.00'
Set-TextValue "E29" '  -0.02%  '
Set-TextValue "D30" '6.19'
Set-TextValue "E30" '  +10.36%  '
Set-TextValue "D31" '1.44'
Set-TextValue "E31" '  +6.82%  '
Set-TextValue "E32" '  +0.25%  '
Set-TextValue "D33" '6.55'
Set-TextValue "E33" '  -0.15%  '
Set-TextValue "D34" '23.64'
Set-TextValue "E34" '  +2.25%  '
Set-TextValue "D35" '0.999'
Set-TextValue "E35" '  +0.14%  '
Set-TextValue "D36" '7.10'
Set-TextValue "E36" '  +6.11%  '
Set-TextValue "D37" '1.52'
Set-TextValue "E37" '  +2.05%  '
Set-TextValue "D38" '162.01'
Set-TextValue "E38" '  +2.28%  '
Set-TextValue "D39" '1.91'
Set-TextValue "E39" '  +1.90%  '
Set-TextValue "D40" '0.0772'
Set-TextValue "E40" '  +0.83%  '
Set-TextValue "D41" '27.42'
Set-TextValue "E41" '  -0.31%  '
Set-TextValue "D42" '2.931.93'
Set-TextValue "E42" '  +0.38%  '
Set-TextValue "D43" '4.57'
Set-TextValue "E43" '  +5.95%  '
Set-TextValue "D44" '42.84'
Set-TextValue "E44" '  +3.40%  '
Set-TextValue "E45" '  -1.54%  '
Set-TextValue "D46" '0.776'
Set-TextValue "E46" '  +1.65%  '
Set-TextValue "D47" '24.11'
Set-TextValue "E47" '  +7.29%  '
Set-TextValue "E48" '  +2.47%  '
Set-TextValue "D49" '0.873'
Set-TextValue "E49" '  +6.82%  '
Set-TextValue "D50" '2.19'
Set-TextValue "E50" '  +12.43%  '
Set-TextValue "E51" '  +3.55%  '
